$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price text looks numeric (e.g. "1.91") must be pre-set to
# Text format, otherwise Excel auto-converts the literal to a real number.
$textCells = @(
    'D5',
    'D6',
    'D7',
    'D8',
    'D9',
    'D10',
    'D12',
    'D13',
    'D15',
    'D18',
    'D20',
    'D22',
    'D23',
    'D24',
    'D25',
    'D27',
    'D28',
    'D29',
    'D31',
    'D32',
    'D33',
    'D34',
    'D35',
    'D36',
    'D37',
    'D38',
    'D39',
    'D40',
    'D41',
    'D42',
    'D43',
    'D44',
    'D45',
    'D46',
    'D47',
    'D48',
    'D49',
    'D50',
    'D51'
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '97.344.31'
$ws.Range('E2').Value = '  +0.69%  '

$ws.Range('D3').Value = '3.704.75'
$ws.Range('E3').Value = '  +0.69%  '

$ws.Range('E4').Value = '  +0.02%  '

$ws.Range('D5').Value = '238.29'
$ws.Range('E5').Value = '  -1.30%  '

$ws.Range('D6').Value = '1.91'
$ws.Range('E6').Value = '  +3.56%  '

$ws.Range('D7').Value = '659.97'
$ws.Range('E7').Value = '  -0.45%  '

$ws.Range('D8').Value = '0.424'
$ws.Range('E8').Value = '  +0.12%  '

$ws.Range('B9').Value = 'Cardano'
$ws.Range('C9').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D9').Value = '1.07'
$ws.Range('E9').Value = '  -1.03%  '

$ws.Range('B10').Value = 'USDC'
$ws.Range('C10').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D10').Value = '0.999'
$ws.Range('E10').Value = '  -0.03%  '

$ws.Range('D11').Value = '3.703.28'
$ws.Range('E11').Value = '  +0.67%  '

$ws.Range('D12').Value = '0.0000311'
$ws.Range('E12').Value = '  +14.98%  '

$ws.Range('D13').Value = '44.30'
$ws.Range('E13').Value = '  -2.68%  '

$ws.Range('E14').Value = '  +1.86%  '

$ws.Range('D15').Value = '6.78'
$ws.Range('E15').Value = '  -2.86%  '

$ws.Range('D16').Value = '4.393.73'
$ws.Range('E16').Value = '  +0.73%  '

$ws.Range('D17').Value = '97.116.64'
$ws.Range('E17').Value = '  +0.79%  '

$ws.Range('D18').Value = '9.17'
$ws.Range('E18').Value = '  +2.97%  '

$ws.Range('D19').Value = '3.712.17'
$ws.Range('E19').Value = '  +0.80%  '

$ws.Range('D20').Value = '13.04'
$ws.Range('E20').Value = '  +0.52%  '

$ws.Range('E21').Value = '  +1.43%  '

$ws.Range('D22').Value = '0.507'
$ws.Range('E22').Value = '  -3.83%  '

$ws.Range('D23').Value = '522.88'
$ws.Range('E23').Value = '  -1.32%  '

$ws.Range('D24').Value = '3.45'
$ws.Range('E24').Value = '  +0.21%  '

$ws.Range('D25').Value = '0.0000213'
$ws.Range('E25').Value = '  +4.57%  '

$ws.Range('E26').Value = '  -0.65%  '

$ws.Range('D27').Value = '101.82'
$ws.Range('E27').Value = '  -0.47%  '

$ws.Range('D28').Value = '0.195'
$ws.Range('E28').Value = '  +14.83%  '

$ws.Range('D29').Value = '13.61'
$ws.Range('E29').Value = '  +4.13%  '

$ws.Range('B30').Value = 'WrappedeETH'
$ws.Range('C30').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D30').Value = '3.903.37'
$ws.Range('E30').Value = '  +0.71%  '

$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').Value = '12.63'
$ws.Range('E31').Value = '  +0.38%  '

$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').Value = '3.06'
$ws.Range('E32').Value = '  -0.31%  '

$ws.Range('B33').Value = 'Dai'
$ws.Range('C33').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D33').Value = '1.00'
$ws.Range('E33').Value = '  +0.08%  '

$ws.Range('B34').Value = 'Cronos'
$ws.Range('C34').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D34').Value = '0.191'
$ws.Range('E34').Value = '  +2.61%  '

$ws.Range('B35').Value = 'Fetch.AI'
$ws.Range('C35').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D35').Value = '1.90'
$ws.Range('E35').Value = '  +2.56%  '

$ws.Range('B36').Value = 'Binance-PegBSC-USD'
$ws.Range('C36').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D36').Value = '1.00'
$ws.Range('E36').Value = '  -0.01%  '

$ws.Range('B37').Value = 'Bittensor'
$ws.Range('C37').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D37').Value = '656.58'
$ws.Range('E37').Value = '  +3.20%  '

$ws.Range('B38').Value = 'EthereumClassic'
$ws.Range('C38').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D38').Value = '32.28'
$ws.Range('E38').Value = '  -1.96%  '

$ws.Range('B39').Value = 'PolygonEcosystemToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D39').Value = '0.597'
$ws.Range('E39').Value = '  +0.30%  '

$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D40').Value = '8.90'
$ws.Range('E40').Value = '  +1.93%  '

$ws.Range('B41').Value = 'USDe'
$ws.Range('C41').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D41').Value = '1.00'
$ws.Range('E41').Value = '  +0.03%  '

$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').Value = '0.165'
$ws.Range('E42').Value = '  +2.19%  '

$ws.Range('D43').Value = '2.06'
$ws.Range('E43').Value = '  +4.45%  '

$ws.Range('B44').Value = 'Filecoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D44').Value = '6.84'
$ws.Range('E44').Value = '  +6.21%  '

$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Value = '40.75'
$ws.Range('E45').Value = '  -9.11%  '

$ws.Range('B46').Value = 'Algorand'
$ws.Range('C46').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D46').Value = '0.491'
$ws.Range('E46').Value = '  +4.62%  '

$ws.Range('B47').Value = 'ARBITRUM'
$ws.Range('C47').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D47').Value = '0.974'
$ws.Range('E47').Value = '  +0.50%  '

$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').Value = '0.0461'
$ws.Range('E48').Value = '  +0.54%  '

$ws.Range('B49').Value = 'Stacks'
$ws.Range('C49').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D49').Value = '2.31'
$ws.Range('E49').Value = '  +0.87%  '

$ws.Range('B50').Value = 'WhiteBITCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D50').Value = '23.66'
$ws.Range('E50').Value = '  +0.03%  '

$ws.Range('B51').Value = 'Cosmos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D51').Value = '8.73'
$ws.Range('E51').Value = '  +1.10%  '

foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
